# "add content MultiLinguals Seeder"
#
# The source workbook uses column B on each sheet ("Blog", "Page",
# "BlogPost") to hold a field-type tag for a list of MultiLingual content
# fields (column A). The original ad-hoc tags (image/int/date/bool/string)
# are replaced with the ABP MultiLingual-seeder property-type names
# (ImageUrl/Text/DateTime/PhoneNumber/EmailAddress), one column-B cell at a
# time, following the semantic meaning of the column-A field name:
#   photo/picture fields  -> ImageUrl
#   phone fields          -> PhoneNumber
#   e-mail fields         -> EmailAddress
#   date fields           -> DateTime
#   everything else       -> Text

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Blog sheet
# ---------------------------------------------------------------------
$wsBlog = $wb.Worksheets.Item("Blog")

$wsBlog.Range("B3").Value = "DateTime"
$wsBlog.Range("B4").Value = "Text"
$wsBlog.Range("B5").Value = "Text"

# ---------------------------------------------------------------------
# Page sheet
# ---------------------------------------------------------------------
$wsPage = $wb.Worksheets.Item("Page")

$wsPage.Range("B3").Value  = "ImageUrl"
$wsPage.Range("B4").Value  = "Text"
$wsPage.Range("B5").Value  = "Text"
$wsPage.Range("B6").Value  = "Text"
$wsPage.Range("B7").Value  = "Text"
$wsPage.Range("B8").Value  = "Text"
$wsPage.Range("B9").Value  = "Text"
$wsPage.Range("B10").Value = "Text"
$wsPage.Range("B11").Value = "PhoneNumber"
$wsPage.Range("B12").Value = "EmailAddress"
$wsPage.Range("B13").Value = "Text"

$wsPage.Range("B15").Value = "ImageUrl"
$wsPage.Range("B16").Value = "Text"
$wsPage.Range("B17").Value = "Text"
$wsPage.Range("B18").Value = "PhoneNumber"
$wsPage.Range("B19").Value = "EmailAddress"

$wsPage.Range("B21").Value = "ImageUrl"
$wsPage.Range("B22").Value = "Text"
$wsPage.Range("B23").Value = "Text"
$wsPage.Range("B24").Value = "Text"
$wsPage.Range("B25").Value = "Text"
$wsPage.Range("B26").Value = "Text"
$wsPage.Range("B27").Value = "PhoneNumber"
$wsPage.Range("B28").Value = "EmailAddress"

$wsPage.Range("B30").Value = "Text"
$wsPage.Range("B31").Value = "Text"
$wsPage.Range("B32").Value = "Text"
$wsPage.Range("B33").Value = "Text"
$wsPage.Range("B34").Value = "Text"
$wsPage.Range("B35").Value = "Text"
$wsPage.Range("B36").Value = "Text"
$wsPage.Range("B37").Value = "Text"

# ---------------------------------------------------------------------
# BlogPost sheet
# ---------------------------------------------------------------------
$wsPost = $wb.Worksheets.Item("BlogPost")

$wsPost.Range("B3").Value = "Text"
$wsPost.Range("B4").Value = "DateTime"
$wsPost.Range("B5").Value = "Text"
$wsPost.Range("B6").Value = "ImageUrl"
$wsPost.Range("B7").Value = "Text"
$wsPost.Range("B8").Value = "Text"

# ---------------------------------------------------------------------
# Page setup for the Blog sheet (Page / BlogPost already carry this)
# ---------------------------------------------------------------------
$wsBlog.PageSetup.PaperSize = 119
$wsBlog.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Restore selections / active cells and make BlogPost the active tab,
# matching the saved UI state in the edited workbook.
# ---------------------------------------------------------------------
$wsBlog.Activate()
$wsBlog.Range("B5").Select()

$wsPage.Activate()
$wsPage.Range("B12").Select()

$wsPost.Activate()
$wsPost.Range("B13").Select()
